$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("consignmentnumber")

$ws.Range("A12").Value = "Reshma"
$ws.Range("F1").Value = "Ayesha"
$ws.Range("G12").Value = "CARDITRECEPTACLEID00000111037"
$ws.Range("G1").Value = "Heena"
$ws.Range("F12").Value = "PREC01048239"

$ws.Columns.Item(6).ColumnWidth = 13.140625
$ws.Columns.Item(7).ColumnWidth = 30.85546875

$ws.Range("G1").Select()
